$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Shared-string-backed header cells on Sheet1 (order matters: controls the
# new shared-string table indices 90/91/92 assigned to "Alpine"/"grasslands"/"plains") ---
$ws1.Range("J1").Value = "Alpine"
$ws1.Range("H1").Value = "grasslands"
$ws1.Range("L1").Value = "plains"

# --- Data edits on Sheet1 row 18 (adds C18, changes K18; C20 COUNTA recalcs automatically) ---
$ws1.Range("C18").Value = 20
$ws1.Range("K18").Value = 20

# --- Column width tweaks on Sheet1 ---
$ws1.Columns.Item(1).ColumnWidth = 17
$ws1.Columns.Item(6).ColumnWidth = 4

# --- View/selection state: Sheet2 loses the active tab, Sheet1 gains it ---
$ws2.Activate() | Out-Null
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 5
$win2.ScrollColumn = 8
$ws2.Range("T4").Select() | Out-Null

$ws1.Activate() | Out-Null
$win1 = $excel.ActiveWindow
$win1.ScrollRow = 2
$win1.ScrollColumn = 2
$ws1.Range("M18").Select() | Out-Null
